# Insert a new weekly record at row 159, pushing the existing rows
# 159-229 down to 160-230 (same behaviour as right-clicking the row
# header and choosing "Insert").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(159).Insert()

# Populate the newly inserted row with the new "Perejil" price record.
$ws.Cells.Item(159, 1).Value = 10
$ws.Cells.Item(159, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(159, 3).Value = "La Araucanía"
$ws.Cells.Item(159, 4).Value = 44523
$ws.Cells.Item(159, 5).Value = 9
$ws.Cells.Item(159, 6).Value = 100112044
$ws.Cells.Item(159, 7).Value = "Perejil"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 20
$ws.Cells.Item(159, 11).Value = 5000
$ws.Cells.Item(159, 12).Value = 5000
$ws.Cells.Item(159, 13).Value = 5000
$ws.Cells.Item(159, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(159, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(159, 16).Value = 1667
$ws.Cells.Item(159, 17).Value = 3
$ws.Cells.Item(159, 18).Value = "Hortaliza"
